$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of row 2 (yellow "input" text style) down into row 3
$ws.Range("A2:H2").Copy()
$ws.Range("A3:H3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New record: same time/detail values as row 2, but a new client code
$ws.Range("A3").Value = "4656709"
$ws.Range("B3").Value = "9"
$ws.Range("C3").Value = "00"
$ws.Range("D3").Value = "12"
$ws.Range("E3").Value = "00"
$ws.Range("F3").Value = "3"
$ws.Range("G3").Value = "00"
$ws.Range("H3").Value = "Conforme"

# Match the author's resulting selection
$ws.Range("B2").Select()
